# Overhaul partial-quotient interface, more interactive lessons
# Rewrites the "en" language sheet Key/Value rows to the new copy, inserting
# the new interactive-lesson rows and re-flowing everything below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C below the header used to hold a few VoiceDuration numbers on
# rows 92-94; those rows move further down the sheet, so clear the whole
# stale range before re-populating it at its new location.
$ws.Range("C2:C103").ClearContents() | Out-Null

$ws.Range("A2").Value2 = 'title'
$ws.Range("B2").Value2 = '<size=50>Attack on Blob</size>\nDivide and Conquer'

$ws.Range("A3").Value2 = 'credits'
$ws.Range("B3").Value2 = 'Made by: RENEGADEWARE'

$ws.Range("A4").Value2 = 'new'
$ws.Range("B4").Value2 = 'NEW GAME'

$ws.Range("A5").Value2 = 'continue'
$ws.Range("B5").Value2 = 'CONTINUE'

$ws.Range("A6").Value2 = 'options'
$ws.Range("B6").Value2 = 'OPTIONS'

$ws.Range("A7").Value2 = 'music'
$ws.Range("B7").Value2 = 'MUSIC'

$ws.Range("A8").Value2 = 'sound'
$ws.Range("B8").Value2 = 'SOUND'

$ws.Range("A9").Value2 = 'speech'
$ws.Range("B9").Value2 = 'SPEECH'

$ws.Range("A10").Value2 = 'close'
$ws.Range("B10").Value2 = 'CLOSE'

$ws.Range("A11").Value2 = 'on'
$ws.Range("B11").Value2 = 'ON'

$ws.Range("A12").Value2 = 'off'
$ws.Range("B12").Value2 = 'OFF'

$ws.Range("A13").Value2 = 'enter'
$ws.Range("B13").Value2 = 'ENTER'

$ws.Range("A14").Value2 = 'cancel'
$ws.Range("B14").Value2 = 'CANCEL'

$ws.Range("A15").Value2 = 'split'
$ws.Range("B15").Value2 = 'SPLIT'

$ws.Range("A16").Value2 = 'victory'
$ws.Range("B16").Value2 = 'VICTORY'

$ws.Range("A17").Value2 = 'attacks'
$ws.Range("B17").Value2 = 'ATTACKS'

$ws.Range("A18").Value2 = 'errors'
$ws.Range("B18").Value2 = 'ERRORS'

$ws.Range("A19").Value2 = 'efficient'
$ws.Range("B19").Value2 = 'EFFICIENT'

$ws.Range("A20").Value2 = 'score'
$ws.Range("B20").Value2 = 'SCORE'

$ws.Range("A21").Value2 = 'total_errors'
$ws.Range("B21").Value2 = 'TOTAL ERROR'

$ws.Range("A22").Value2 = 'total_score'
$ws.Range("B22").Value2 = 'TOTAL SCORE'

$ws.Range("A23").Value2 = 'health_warning'
$ws.Range("B23").Value2 = 'Watch out! Once the health bar is empty, you will have to start over!'

$ws.Range("A24").Value2 = 'placeValue'
$ws.Range("B24").Value2 = 'Place Value'

$ws.Range("A25").Value2 = 'placeValueDistribute'
$ws.Range("B25").Value2 = 'Place Value Distribute'

$ws.Range("A26").Value2 = 'placeOnes'
$ws.Range("B26").Value2 = 'Ones'

$ws.Range("A27").Value2 = 'placeTens'
$ws.Range("B27").Value2 = 'Tens'

$ws.Range("A28").Value2 = 'placeHundreds'
$ws.Range("B28").Value2 = 'Hundreds'

$ws.Range("A29").Value2 = 'areaModel'
$ws.Range("B29").Value2 = 'Area Model'

$ws.Range("A30").Value2 = 'areaModelDistribute'
$ws.Range("B30").Value2 = 'Area Model Distribute'

$ws.Range("A31").Value2 = 'opDivision'
$ws.Range("B31").Value2 = 'Operation Divide'

$ws.Range("A32").Value2 = 'dividend'
$ws.Range("B32").Value2 = 'Dividend'

$ws.Range("A33").Value2 = 'divisor'
$ws.Range("B33").Value2 = 'Divisor'

$ws.Range("A34").Value2 = 'quotient'
$ws.Range("B34").Value2 = 'Quotient'

$ws.Range("A35").Value2 = 'intro_alert_0'
$ws.Range("B35").Value2 = 'Danger!'

$ws.Range("A36").Value2 = 'intro_alert_1'
$ws.Range("B36").Value2 = 'Multiple blobs of epic proportion detected!'

$ws.Range("A37").Value2 = 'intro_alert_2'
$ws.Range("B37").Value2 = 'It seems these blobs came from a lab and have somehow grown uncontrollably.'

$ws.Range("A38").Value2 = 'intro_alert_3'
$ws.Range("B38").Value2 = 'We must find a way to shrink them back!'

$ws.Range("A39").Value2 = 'intro_op_0'
$ws.Range("B39").Value2 = 'According to the latest studies of blob-ology, the only way to shrink them is by the power of the division operation.'

$ws.Range("A40").Value2 = 'intro_op_1'
$ws.Range("B40").Value2 = 'Divide and conquer, as they say!'

$ws.Range("A41").Value2 = 'intro_op_2'
$ws.Range("B41").Value2 = 'Go forth, our intrepid hero! Use your mathematical might to vanquish these menacing blobs!'

$ws.Range("A42").Value2 = 'lesson1_intro_0'
$ws.Range("B42").Value2 = 'Dividing a large number can be daunting, but with the right trick, it can be a breeze!'

$ws.Range("A43").Value2 = 'lesson1_intro_1'
$ws.Range("B43").Value2 = 'Let’s briefly examine a useful mathematical principle to help us.'

$ws.Range("A44").Value2 = 'lesson1_placeValue_0'
$ws.Range("B44").Value2 = 'First, let’s take a quick look at how large numbers are arranged by single digit numbers.'

$ws.Range("A45").Value2 = 'lesson1_placeValue_1'
$ws.Range("B45").Value2 = 'These digits are placed by multiples of 10, since we use a base 10 number system.'

$ws.Range("A46").Value2 = 'lesson1_placeValueDist_0'
$ws.Range("B46").Value2 = 'Let’s put this into practice by splitting up the large number. Drag the equation all the way to the left.'

$ws.Range("A47").Value2 = 'lesson1_digitSwapFirst_0'
$ws.Range("B47").Value2 = 'Now we can start moving the digits from one number into another.'

$ws.Range("A48").Value2 = 'lesson1_digitSwapFirst_1'
$ws.Range("B48").Value2 = 'Press the highlighted digit to make the move.'

$ws.Range("A49").Value2 = 'lesson1_digitSwapSecond_0'
$ws.Range("B49").Value2 = 'Now for the next digit, moving this will still make both numbers wholly divisible.'

$ws.Range("A50").Value2 = 'lesson1_digitSwapComplete_0'
$ws.Range("B50").Value2 = 'As you can see, we now have two divisions that are much easier to solve.'

$ws.Range("A51").Value2 = 'lesson1_digitSwapComplete_1'
$ws.Range("B51").Value2 = 'Press each one to solve the division.'

$ws.Range("A52").Value2 = 'lesson1_divisionsSolved_0'
$ws.Range("B52").Value2 = 'Now we just have to add both numbers to get the final answer.'

$ws.Range("A53").Value2 = 'lesson1_divisionsSolved_1'
$ws.Range("B53").Value2 = 'Why don’t you do the honors by pressing on the plus sign.'

$ws.Range("A54").Value2 = 'lesson1_addSolved_0'
$ws.Range("B54").Value2 = 'Not bad! As you can see, splitting up a large number this way can help solve divisions easily.'

$ws.Range("A55").Value2 = 'lesson1_addSolved_1'
$ws.Range("B55").Value2 = 'Let’s go ahead and put this into practice when we face the mega blob!'

$ws.Range("A56").Value2 = 'lesson2_intro_0'
$ws.Range("B56").Value2 = 'This time around, we will be dividing with double-digit divisors.'

$ws.Range("A57").Value2 = 'lesson2_intro_1'
$ws.Range("B57").Value2 = 'It would be too troublesome to deal with these blobs using our current technique.'

$ws.Range("A58").Value2 = 'lesson2_intro_2'
$ws.Range("B58").Value2 = 'Fortunately, we have one more trick up our sleeves!'

$ws.Range("A59").Value2 = 'lesson2_areaModel_0'
$ws.Range("B59").Value2 = 'Since division is the inverse of multiplication, you can visualize the equation as the dimensions of an area.'

$ws.Range("A60").Value2 = 'lesson2_areaModel_1'
$ws.Range("B60").Value2 = 'In this case, the quotient of the division is the width of the area.'

$ws.Range("A61").Value2 = 'lesson2_areaModel_drag_0'
$ws.Range("B61").Value2 = 'We can use this model to partially solve the division with smaller numbers.'

$ws.Range("A62").Value2 = 'lesson2_areaModel_drag_1'
$ws.Range("B62").Value2 = 'Drag the area from left to right to see how this works.'

$ws.Range("A63").Value2 = 'lesson2_areaModel_drag_complete_0'
$ws.Range("B63").Value2 = 'As you can see, each number multiplied by 12, the divisor, is subtracted from the dividend.'

$ws.Range("A64").Value2 = 'lesson2_areaModel_drag_complete_1'
$ws.Range("B64").Value2 = 'Now we can easily divide the remaining dividend to get the final answer.'

$ws.Range("A65").Value2 = 'lesson2_areaModel_answer_0'
$ws.Range("B65").Value2 = 'Adding the split values will then give you the whole answer.'

$ws.Range("A66").Value2 = 'lesson2_end_0'
$ws.Range("B66").Value2 = 'Now why don’t we try this new technique with the next mega blob!'

$ws.Range("A67").Value2 = 'level1_intro_0'
$ws.Range("B67").Value2 = 'Look out! Two blobs have appeared.'

$ws.Range("A68").Value2 = 'level1_intro_1'
$ws.Range("B68").Value2 = 'In order to beat the mega blob, we must merge all the blobs into one final quotient blob.'

$ws.Range("A69").Value2 = 'level1_intro_2'
$ws.Range("B69").Value2 = 'Those numbers don’t look that scary. We can directly solve the division.'

$ws.Range("A70").Value2 = 'instruct_drag_blob_0'
$ws.Range("B70").Value2 = 'In order to merge the two blobs, simply drag one to another like so.'

$ws.Range("A71").Value2 = 'op_instruct_0'
$ws.Range("B71").Value2 = 'Now you must solve the operation by typing in the number via the numpad.'

$ws.Range("A72").Value2 = 'op_instruct_1'
$ws.Range("B72").Value2 = 'You can also use the keyboard to enter the numbers.'

$ws.Range("A73").Value2 = 'op_instruct_2'
$ws.Range("B73").Value2 = 'Once you feel confident with your answer, press the ENTER button on the numpad (or your keyboard).'

$ws.Range("A74").Value2 = 'attack_instruct_success_0'
$ws.Range("B74").Value2 = 'Excellent! We’ve managed to clear the blobs!'

$ws.Range("A75").Value2 = 'attack_instruct_boss_hp_0'
$ws.Range("B75").Value2 = 'This is the representation of the mega blob’s health.'

$ws.Range("A76").Value2 = 'attack_instruct_boss_hp_1'
$ws.Range("B76").Value2 = 'As you can see, it has been reduced.'

$ws.Range("A77").Value2 = 'attack_instruct_boss_hp_2'
$ws.Range("B77").Value2 = 'Once it’s empty, the mega blob will be defeated.'

$ws.Range("A78").Value2 = 'split_instruct_0'
$ws.Range("B78").Value2 = 'Now we are dealing with a much larger blob!'

$ws.Range("A79").Value2 = 'split_instruct_1'
$ws.Range("B79").Value2 = 'Let’s split the blob up into two to make our life easier.'

$ws.Range("A80").Value2 = 'split_instruct_2'
$ws.Range("B80").Value2 = 'Press the sparkly blob as shown to proceed.'

$ws.Range("A81").Value2 = 'split_op_instruct_0'
$ws.Range("B81").Value2 = 'Here you can see a representation of how the blob is going to be split.'

$ws.Range("A82").Value2 = 'split_op_instruct_1'
$ws.Range("B82").Value2 = 'We will be splitting the blob by transferring its digits to a new blob.'

$ws.Range("A83").Value2 = 'split_op_instruct_2'
$ws.Range("B83").Value2 = 'Simply click on any of the digits to transfer them.'

$ws.Range("A84").Value2 = 'split_op_instruct_3'
$ws.Range("B84").Value2 = 'Once you are happy with the new split numbers, press the SPLIT button.'

$ws.Range("A85").Value2 = 'split_op_instruct_4'
$ws.Range("B85").Value2 = 'Remember, both new numbers must be wholly divisible for the split to succeed!'

$ws.Range("A86").Value2 = 'split_op_success_0'
$ws.Range("B86").Value2 = 'Nicely done! The blobs have now been split into two.'

$ws.Range("A87").Value2 = 'split_op_success_1'
$ws.Range("B87").Value2 = 'You can split the blobs further if you want, but there’s a limit!'

$ws.Range("A88").Value2 = 'split_op_success_2'
$ws.Range("B88").Value2 = 'Once all the blobs have been merged into the final quotient blob, we will be able to attack.'

$ws.Range("A89").Value2 = 'split_op_success_3'
$ws.Range("B89").Value2 = 'Good luck!'

$ws.Range("A90").Value2 = 'level3_intro_0'
$ws.Range("B90").Value2 = 'Watch out! These blobs are not to be trifled with!'

$ws.Range("A91").Value2 = 'level3_intro_1'
$ws.Range("B91").Value2 = 'We’ll be employing the area model trick we just learned to defeat this blob.'

$ws.Range("A92").Value2 = 'split_instruct2_0'
$ws.Range("B92").Value2 = 'Just as you have done before, press the sparkly blob to commence the split.'

$ws.Range("A93").Value2 = 'split_op_partial_instruct_0'
$ws.Range("B93").Value2 = 'Now we are going to reduce the dividend blob’s number by multiplying the divisor blob with a number.'

$ws.Range("A94").Value2 = 'split_op_partial_instruct_1'
$ws.Range("B94").Value2 = 'The best approach is to multiply a single number by 10 several times, as long as it’s not larger than the dividend.'

$ws.Range("A95").Value2 = 'split_op_partial_instruct_2'
$ws.Range("B95").Value2 = 'We will do exactly that for this problem!'

$ws.Range("A96").Value2 = 'split_op_partial_mult_tens'
$ws.Range("B96").Value2 = 'Press the left button until the number is 100.'

$ws.Range("A97").Value2 = 'split_op_partial_mult_digit'
$ws.Range("B97").Value2 = 'Now press the up button until the number is 200.'

$ws.Range("A98").Value2 = 'split_op_partial_next'
$ws.Range("B98").Value2 = 'Press this button to proceed.'

$ws.Range("A99").Value2 = 'split_op_partial_next_instruct_0'
$ws.Range("B99").Value2 = 'Now you must type in the correct answer for the multiplication.'

$ws.Range("A100").Value2 = 'split_op_partial_next_instruct_1'
$ws.Range("B100").Value2 = 'Once you press ENTER, you’ll see that number subtracted from the dividend blob.'

$ws.Range("A101").Value2 = 'split_op_partial_success_0'
$ws.Range("B101").Value2 = 'Good! The blob’s value has been reduced, and a partial quotient blob has appeared.'

$ws.Range("A102").Value2 = 'split_op_partial_success_1'
$ws.Range("B102").Value2 = 'Splitting up large numbers this way will make it easier to deal with two or more-digit divisors.'

$ws.Range("A103").Value2 = 'split_op_partial_success_2'
$ws.Range("B103").Value2 = 'Remember this trick, and you will be able to defeat these mega blobs with ease!'

$ws.Range("A104").Value2 = 'end_congrats'
$ws.Range("B104").Value2 = 'CONGRATULATIONS!'
$ws.Range("C104").Value2 = 2

$ws.Range("A105").Value2 = 'end_desc'
$ws.Range("B105").Value2 = 'You have vanquished all the mega blobs! Peace has returned to the world once more!'
$ws.Range("C105").Value2 = 5

$ws.Range("A106").Value2 = 'end_thanks'
$ws.Range("B106").Value2 = 'Thank you for playing!'
$ws.Range("C106").Value2 = 2

# Match the author's final cursor position/selection on the sheet.
$ws.Range("B103").Select() | Out-Null
